# outputs-r202/test-g__RUG033.xlsx was refreshed by the batch pipeline: the
# "quadratic-svm-score" sheet gets re-stamped with a fresh style pass (the
# text-format header/row-label cells pick up a new cell style further along
# the shared style table) and the predicted score in B2 is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text number format to the header row and the row-label cell
# so they pick up a newly minted (but format-equivalent) text style, just as
# the refreshed pipeline run re-stamped them.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2").NumberFormat = "@"

# The refreshed prediction score for RUG033.
$ws.Range("B2").Value = 0.90730389906653741
